$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), copying the header style (bold/border)
# used by the rest of row 1 (e.g. H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J for rows 2-9
$values = @{
    2 = @{ I = 6; J = 8 }
    3 = @{ I = 9; J = 9 }
    4 = @{ I = 8; J = 8 }
    5 = @{ I = 8; J = 9 }
    6 = @{ I = 8; J = 8 }
    7 = @{ I = 9; J = 9 }
    8 = @{ I = 8; J = 8 }
    9 = @{ I = 6; J = 6 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 9).Value = $values[$row].I
    $ws.Cells.Item($row, 10).Value = $values[$row].J
}
